$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into the cell at $addr while forcing text storage.
#
# Every Price / Volume(1h) cell on this sheet is stored as text, even when
# the text happens to look like a plain number (e.g. "58.70", "1.00",
# "0.999"). Assigning such a string straight to Range.Value lets Excel
# auto-coerce it to a Double, silently destroying information that matters
# here (58.70 -> 58.7, 1.00 -> 1). Prefixing the literal with a leading
# apostrophe is the normal Excel convention for forcing a numeric-looking
# value to stay text - exactly as if it had been typed into the UI that way
# - so we do that whenever the text looks like a plain decimal number.
function Set-CellText {
    param([string]$addr, [string]$text)
    if ($text -match '^-?\d+(\.\d+)?$') {
        $ws.Range($addr).Value = "'" + $text
    } else {
        $ws.Range($addr).Value = $text
    }
}

Set-CellText "D2" '69.071.97'
Set-CellText "E2" '  +5.24%  '
Set-CellText "D3" '3.543.04'
Set-CellText "E3" '  +4.22%  '
Set-CellText "D4" '0.999'
Set-CellText "E4" '  -0.07%  '
Set-CellText "D5" '590.09'
Set-CellText "E5" '  +5.24%  '
Set-CellText "D6" '192.44'
Set-CellText "E6" '  +9.07%  '
Set-CellText "E7" '  +1.07%  '
Set-CellText "D8" '3.531.53'
Set-CellText "E8" '  +4.14%  '
Set-CellText "D9" '0.999'
Set-CellText "E9" '  -0.09%  '
Set-CellText "D10" '0.178'
Set-CellText "E10" '  +3.53%  '
Set-CellText "D11" '0.658'
Set-CellText "E11" '  +2.69%  '
Set-CellText "D12" '58.70'
Set-CellText "E12" '  +9.27%  '
Set-CellText "D13" '0.0000292'
Set-CellText "E13" '  +4.97%  '
Set-CellText "D14" '9.63'
Set-CellText "E14" '  +4.44%  '
Set-CellText "D15" '4.099.39'
Set-CellText "E15" '  +3.89%  '
Set-CellText "D16" '19.18'
Set-CellText "E16" '  +4.54%  '
Set-CellText "D17" '3.540.82'
Set-CellText "E17" '  +3.79%  '
Set-CellText "D18" '69.043.84'
Set-CellText "E18" '  +5.34%  '
Set-CellText "D19" '12.38'
Set-CellText "E19" '  +4.25%  '
Set-CellText "E20" '  +0.29%  '
Set-CellText "D22" '493.95'
Set-CellText "E22" '  +2.46%  '
Set-CellText "D23" '5.65'
Set-CellText "E23" '  +13.95%  '
Set-CellText "D24" '17.41'
Set-CellText "E24" '  +21.60%  '
Set-CellText "D25" '4.44'
Set-CellText "E25" '  +7.74%  '
Set-CellText "D26" '91.03'
Set-CellText "E26" '  +2.12%  '
Set-CellText "D27" '3.04'
Set-CellText "E27" '  +4.07%  '
Set-CellText "D28" '11.13'
Set-CellText "D29" '9.25'
Set-CellText "E29" '  +5.75%  '
Set-CellText "D30" '31.83'
Set-CellText "E30" '  +1.10%  '
Set-CellText "D31" '7.47'
Set-CellText "E31" '  +13.86%  '
Set-CellText "D32" '610.68'
Set-CellText "E32" '  +5.99%  '
Set-CellText "D33" '12.06'
Set-CellText "E33" '  +4.55%  '
Set-CellText "D34" '65.20'
Set-CellText "E34" '  +4.38%  '
Set-CellText "E35" '  +5.37%  '
Set-CellText "E36" '  +5.23%  '
Set-CellText "E37" '  +0.00%  '
Set-CellText "D38" '37.70'
Set-CellText "E38" '  +4.59%  '
Set-CellText "D39" '0.0₃0796'
Set-CellText "E39" '  +7.71%  '
Set-CellText "E40" '  +4.98%  '
Set-CellText "E41" '  -1.41%  '
Set-CellText "D42" '3.278.13'
Set-CellText "E42" '  +5.10%  '
Set-CellText "D43" '3.00'
Set-CellText "E43" '  +7.72%  '
Set-CellText "E44" '  +5.27%  '
Set-CellText "D45" '2.67'
Set-CellText "E45" '  +9.05%  '
Set-CellText "D46" '3.30'
Set-CellText "E46" '  +4.03%  '
Set-CellText "E47" '  +1.71%  '
Set-CellText "D48" '2.75'
Set-CellText "E48" '  +17.91%  '
Set-CellText "D49" '9.06'
Set-CellText "E49" '  +7.29%  '
Set-CellText "D50" '1.00'
Set-CellText "E50" '  +0.26%  '
Set-CellText "D51" '141.04'
Set-CellText "E51" '  +0.50%  '
